$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = [double]"0.9999507296423068"
$ws.Range("E3").Value = [double]"0.9999507296423068"

$ws.Range("D4").Value = [double]"0.999976429191315"
$ws.Range("E4").Value = [double]"0.999976429191315"

$ws.Range("D5").Value = [double]"3.239886489345498E-09"
$ws.Range("E5").Value = [double]"3.239886489345498E-09"

$ws.Range("D6").Value = [double]"1.256613515613315E-10"
$ws.Range("E6").Value = [double]"1.256613515613315E-10"

$ws.Range("D7").Value = [double]"0.0005057763194410223"
$ws.Range("E7").Value = [double]"0.999494223680559"

$ws.Range("D8").Value = [double]"0.9999999996400917"
$ws.Range("E8").Value = [double]"3.599083253646995E-10"

$ws.Range("D9").Value = [double]"0.9999999986786048"
$ws.Range("E9").Value = [double]"1.321395215470034E-09"

$ws.Range("D10").Value = [double]"3.834142887032388E-05"
$ws.Range("E10").Value = [double]"0.9999616585711297"

$ws.Range("D11").Value = [double]"0.9999920173804305"
$ws.Range("E11").Value = [double]"7.982619569468774E-06"
$ws.Range("F11").Value = [double]"8.070170402526855"
